$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 731
$ws1.Range("F3").Value = 600
$ws1.Range("F4").Value = 559
$ws1.Range("F7").Value = 87
$ws1.Range("F9").Value = 55
$ws1.Range("F10").Value = 8
$ws1.Range("F11").Value = 4841
$ws1.Range("F12").Value = 4566
$ws1.Range("F13").Value = 7
$ws1.Range("F14").Value = 24
$ws1.Range("F16").Value = 36
$ws1.Range("F17").Value = 167

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 67

# --- Sheet "全部类型" (all types, combined) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 731
$ws4.Range("F3").Value = 600
$ws4.Range("F4").Value = 559
$ws4.Range("F7").Value = 87
$ws4.Range("F9").Value = 55
$ws4.Range("F10").Value = 8
$ws4.Range("F11").Value = 4841
$ws4.Range("F12").Value = 4566
$ws4.Range("F13").Value = 7
$ws4.Range("F14").Value = 24
$ws4.Range("F16").Value = 36
$ws4.Range("F17").Value = 167
$ws4.Range("F18").Value = 67
